$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Access")
$ws.Range("C2").Value = "Zarif@mohd"
$ws.Hyperlinks.Add($ws.Range("C2"), "http://rmstest.ehealthcorp.net:8020/Default.aspx")
$ws.Range("C2").Style = "Hyperlink"
